# "Version 2." -> "Version 1."  (revert of "Wireframes version 2")
#
# The original paragraph is built from these runs:
#   "Versi" | "on" | " 2" | "." (plus a _GoBack bookmark sitting between
#   the " 2" run and the "." run)
# The target paragraph collapses "Versi"+"on" into a single "Version" run,
# turns " 2" into " 1." and removes the trailing "." run altogether, so the
# paragraph ends right after the bookmark.

$d = $word.ActiveDocument

# 1) Merge the "Versi" + "on" runs into a single "Version" run.
#    Setting identical text is a no-op here, so nudge it through a
#    temporary value first, forcing the two runs to coalesce into one.
$r = $d.Range(0, 7)
$r.Text = "Versionx"
$r = $d.Range(0, 8)
$r.Text = "Version"

# 2) Turn the "2" in " 2" into "1." (yields "Version 1..", i.e. the
#    trailing "." run is still present after this step).
$r = $d.Range(8, 9)
$r.Text = "1."

# 3) Delete the now-redundant trailing "." run (the extra period that
#    used to follow the bookmark).
$r = $d.Range(10, 11)
$r.Delete()
